$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the three data rows (frances.burns / joe.larson / cedric.kelly)
# into rows 5-7 below the existing rows 2-4
$ws.Range("A5").Value = $ws.Range("A2").Value2
$ws.Range("B5").Value = $ws.Range("B2").Value2
$ws.Range("A6").Value = $ws.Range("A3").Value2
$ws.Range("B6").Value = $ws.Range("B3").Value2
$ws.Range("A7").Value = $ws.Range("A4").Value2
$ws.Range("B7").Value = $ws.Range("B4").Value2

# New "status"/"role" columns (C/D/E) for every data row are quote-prefixed
# empty-text cells (adds the blank shared string before the new header text)
$ws.Range("C2:E7").Value = "'"

# New header cells: status / role, styled like the existing header row
$ws.Range("D1").Value = "status"
$ws.Range("E1").Value = "role"
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6").Select()
